$d = $word.ActiveDocument

function Get-ParaByPrefix($prefix) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($prefix)) {
            return $p
        }
    }
    return $null
}

function Set-ParaText($prefix, $newText) {
    $p = Get-ParaByPrefix($prefix)
    if ($p -eq $null) {
        throw "Paragraph not found for prefix: $prefix"
    }
    $p.Range.Text = $newText
}

# 1) Update the "Generated" timestamp.
Set-ParaText 'Generated: 2026-02-18 23:49 UTC' 'Generated: 2026-02-19 01:32 UTC'

# 2) Update the confidence / veracity score.
Set-ParaText 'Confidence / veracity: 75/100' 'Confidence / veracity: 52/100'

# 3) Swap the order of the "Uber Liability Verdict" and "Consumer Watchdog's
#    Consumer Alert" headline paragraphs.
$pLiability = Get-ParaByPrefix('- Uber Liability Verdict Puts Safety Data')
$pWatchdog  = Get-ParaByPrefix("- Consumer Watchdog")

$liabilityText = '- Uber Liability Verdict Puts Safety Data And Investor Risks In Focus | finnhub | INSURANCE | impact -2 | https://finnhub.io/api/news?id=f10e3a67099b57552a9d829bee868929c2cc772f1f637964f1cd149e473b8a20'
$watchdogText  = '- Consumer Watchdog''s Consumer Alert Challenges Uber''s Misleading Super Bowl "Consumer Alert" TV Ad: Uber Really Wants Removal of Legal Accountability for Faulty Safety Checks | finnhub | SAFETY | impact -2 | https://finnhub.io/api/news?id=e8d0314a28db6d2e4ebee541d728e8e6e3d7b317a011de45ec5be7c443fdc36b'

$pLiability.Range.Text = $watchdogText
$pWatchdog.Range.Text  = $liabilityText

# 4) Replace the five "OTHER/FINANCIAL" headline paragraphs with the new set
#    of headlines.
Set-ParaText '- David Tepper' '- IBM spin-off CFO departs amid cash management review: Trial Balance | finnhub | OTHER | impact 0 | https://finnhub.io/api/news?id=82270a4c740bf3b7678322f439605063222c3186cb41c0c3f3796392b0221292'

Set-ParaText "- Uber offers incentives for EV chargers" ('- Uber Deepens T' + [char]0x00FC + 'rkiye Footprint With Getir Delivery Buyout | finnhub | OTHER | impact 0 | https://finnhub.io/api/news?id=33d2c5ad35b4d45f75ed416b52f6132b4b1a2216f7b1fc4b392be4492c910a70')

Set-ParaText '- Uber to Invest' ('- Uber To Acquire Getir''s T' + [char]0x00FC + 'rkiye Delivery Portfolio Covering Food, Grocery, Retail, And Water Services; Terms Not Disclosed | finnhub | OTHER | impact 0 | https://finnhub.io/api/news?id=049ae8617ee0c75685d83074cae5958e3cad04be216687f8b7f0e75e6802259f')

Set-ParaText '- Jim Cramer Discusses' '- Lyft debuts teen accounts more than two years after Uber | finnhub | OTHER | impact 0 | https://finnhub.io/api/news?id=de2d8eae90608556df6105b711a5c099a8ef9319f40596da163d78ed0b36599b'

Set-ParaText '- DoorDash Earnings' '- Tesla Announces New Semi Truck Trim Levels, Elon Musk Reaffirms Volume Production This Year | finnhub | OTHER | impact 0 | https://finnhub.io/api/news?id=cbde98592fd6bccdfe0b7dcae8899aa7b050ceddbb3674619d857471834d7ec4'
